# Generate Report for Archive
# Update the Status for the "b814fa4b-2a4a-4bb4-a9c7-189647d9a0e4.md" row
# (row 6) from "Ready for handoff" to "In Translation" across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B6").Value = "In Translation"
$overview.Range("C6").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B6").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B6").Value = "In Translation"
